# Applies the Betfair Back/Lay odds update for 2025-12-16
# (values updated as odds moved closer to game time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 3.55
$ws.Range("H2").Value = 2.34
$ws.Range("I2").Value = 2.7
$ws.Range("L2").Value = 1.45
$ws.Range("N2").Value = 3.3
$ws.Range("P2").Value = 1.76
$ws.Range("Q2").Value = 2.06
$ws.Range("S2").Value = 3.7
$ws.Range("U2").Value = 1.99
$ws.Range("V2").Value = 1.59
$ws.Range("W2").Value = 1.4
$ws.Range("AF2").Value = 24
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 19.5
$ws.Range("AL2").Value = 150

# Row 3
$ws.Range("F3").Value = 1.48
$ws.Range("G3").Value = 1.56
$ws.Range("H3").Value = 6.8
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 5.2
$ws.Range("L3").Value = 1.41
$ws.Range("N3").Value = 3.6
$ws.Range("O3").Value = 1.31
$ws.Range("P3").Value = 1.89
$ws.Range("Q3").Value = 1.94
$ws.Range("S3").Value = 3.3
$ws.Range("T3").Value = 1.96
$ws.Range("W3").Value = 2.78
$ws.Range("AB3").Value = 16.5
$ws.Range("AC3").Value = 20
$ws.Range("AN3").Value = 29

# Row 4
$ws.Range("F4").Value = 1.27
$ws.Range("G4").Value = 1.34
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 16.5
$ws.Range("J4").Value = 5.4
$ws.Range("K4").Value = 6.8
$ws.Range("L4").Value = 1.34
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 4.1
$ws.Range("P4").Value = 2.12
$ws.Range("R4").Value = 1.43
$ws.Range("S4").Value = 2.8
$ws.Range("T4").Value = 2.22
$ws.Range("U4").Value = 1.62
$ws.Range("AC4").Value = 14.5
$ws.Range("AK4").Value = 16.5
$ws.Range("AN4").Value = 5.9

# Row 5
$ws.Range("G5").Value = 1.66
$ws.Range("L5").Value = 1.42
$ws.Range("N5").Value = 3.35
$ws.Range("P5").Value = 1.84
$ws.Range("Q5").Value = 1.98
$ws.Range("S5").Value = 3.5
$ws.Range("T5").Value = 2
$ws.Range("W5").Value = 2.5

# Row 6
$ws.Range("F6").Value = 2.78
$ws.Range("G6").Value = 3.05
$ws.Range("H6").Value = 2.42
$ws.Range("I6").Value = 2.64
$ws.Range("J6").Value = 3.65
$ws.Range("L6").Value = 1.37
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 4.4
$ws.Range("O6").Value = 1.26
$ws.Range("Q6").Value = 1.79
$ws.Range("R6").Value = 1.45
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 1.69
$ws.Range("V6").Value = 1.61
$ws.Range("W6").Value = 1.5
$ws.Range("X6").Value = 20
$ws.Range("Z6").Value = 18.5
$ws.Range("AA6").Value = 60
$ws.Range("AC6").Value = 8.6
$ws.Range("AE6").Value = 27
$ws.Range("AF6").Value = 22
$ws.Range("AH6").Value = 16.5
$ws.Range("AI6").Value = 60
$ws.Range("AJ6").Value = 130
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 65
$ws.Range("AN6").Value = 25
$ws.Range("AO6").Value = 18.5

# Row 7
$ws.Range("F7").Value = 6.4
$ws.Range("H7").Value = 1.18
$ws.Range("I7").Value = 1.24
$ws.Range("J7").Value = 5.1
$ws.Range("L7").Value = 1.2
$ws.Range("N7").Value = 2.68
$ws.Range("P7").Value = 1.67
$ws.Range("Q7").Value = 1.36
$ws.Range("R7").Value = 1.67
$ws.Range("T7").Value = 1.91
$ws.Range("U7").Value = 1.64
$ws.Range("V7").Value = 3.15
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 14.5
$ws.Range("Z7").Value = 11
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 15.5
$ws.Range("AE7").Value = 17
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 46
$ws.Range("AI7").Value = 48
$ws.Range("AO7").Value = 3.5

# Row 8
$ws.Range("G8").Value = 3.25
$ws.Range("H8").Value = 2.56
$ws.Range("Q8").Value = 2.08
$ws.Range("S8").Value = 3.55
$ws.Range("W8").Value = 1.44
$ws.Range("AA8").Value = 130
$ws.Range("AJ8").Value = 130
$ws.Range("AK8").Value = 100
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 70

# Row 9
$ws.Range("H9").Value = 4.6
$ws.Range("J9").Value = 3.7
$ws.Range("K9").Value = 3.95
$ws.Range("L9").Value = 1.44
$ws.Range("N9").Value = 3.7
$ws.Range("O9").Value = 1.34
$ws.Range("P9").Value = 1.89
$ws.Range("Q9").Value = 2.02
$ws.Range("R9").Value = 1.34
$ws.Range("S9").Value = 3.6
$ws.Range("T9").Value = 1.88
$ws.Range("U9").Value = 1.99
$ws.Range("X9").Value = 14
$ws.Range("Y9").Value = 16
$ws.Range("Z9").Value = 38
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 8.4
$ws.Range("AD9").Value = 19.5
$ws.Range("AH9").Value = 20
$ws.Range("AK9").Value = 21
$ws.Range("AN9").Value = 14
$ws.Range("AO9").Value = 280

# Row 10
$ws.Range("H10").Value = 2.94
$ws.Range("I10").Value = 3.1
$ws.Range("K10").Value = 3.85
$ws.Range("L10").Value = 1.36
$ws.Range("M10").Value = 1.05
$ws.Range("N10").Value = 4.8
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 2.24
$ws.Range("Q10").Value = 1.76
$ws.Range("R10").Value = 1.5
$ws.Range("S10").Value = 2.92
$ws.Range("T10").Value = 1.6
$ws.Range("U10").Value = 2.44
$ws.Range("X10").Value = 19.5
$ws.Range("Y10").Value = 15.5
$ws.Range("Z10").Value = 23
$ws.Range("AB10").Value = 14
$ws.Range("AE10").Value = 34
$ws.Range("AF10").Value = 18
$ws.Range("AH10").Value = 17.5
$ws.Range("AI10").Value = 36
$ws.Range("AJ10").Value = 36
$ws.Range("AK10").Value = 28
$ws.Range("AM10").Value = 70
$ws.Range("AN10").Value = 16.5
$ws.Range("AO10").Value = 23

# Row 11
$ws.Range("H11").Value = 1.91
$ws.Range("K11").Value = 4.4
$ws.Range("L11").Value = 1.29
$ws.Range("O11").Value = 1.18
$ws.Range("Q11").Value = 1.57
$ws.Range("R11").Value = 1.69
$ws.Range("S11").Value = 2.4
$ws.Range("T11").Value = 1.56
$ws.Range("Z11").Value = 15
$ws.Range("AB11").Value = 22
$ws.Range("AC11").Value = 10
$ws.Range("AF11").Value = 34
$ws.Range("AH11").Value = 15.5
$ws.Range("AI11").Value = 26
$ws.Range("AJ11").Value = 85
$ws.Range("AL11").Value = 40
$ws.Range("AN11").Value = 27

# Row 12
$ws.Range("F12").Value = 2.88
$ws.Range("G12").Value = 3.1
$ws.Range("I12").Value = 2.66
$ws.Range("K12").Value = 3.65
$ws.Range("L12").Value = 1.43
$ws.Range("N12").Value = 3.6
$ws.Range("P12").Value = 1.86
$ws.Range("Q12").Value = 2.04
$ws.Range("R12").Value = 1.33
$ws.Range("T12").Value = 1.8
$ws.Range("V12").Value = 1.6
$ws.Range("W12").Value = 1.47
$ws.Range("X12").Value = 15
$ws.Range("AA12").Value = 48
$ws.Range("AH12").Value = 20

# Row 13
$ws.Range("G13").Value = 1.91
$ws.Range("H13").Value = 4.3
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 3.85
$ws.Range("K13").Value = 4.5
$ws.Range("L13").Value = 1.36
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 4.3
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 2.18
$ws.Range("Q13").Value = 1.76
$ws.Range("R13").Value = 1.43
$ws.Range("S13").Value = 2.94
$ws.Range("T13").Value = 1.73
$ws.Range("U13").Value = 2.12
$ws.Range("V13").Value = 1.26
$ws.Range("X13").Value = 19
$ws.Range("Y13").Value = 20
$ws.Range("AB13").Value = 11
$ws.Range("AC13").Value = 9.800000000000001
$ws.Range("AF13").Value = 13
$ws.Range("AH13").Value = 19
$ws.Range("AI13").Value = 160
$ws.Range("AM13").Value = 1000
$ws.Range("AN13").Value = 10.5

# Row 14
$ws.Range("F14").Value = 1.89
$ws.Range("G14").Value = 1.99
$ws.Range("L14").Value = 1.46
$ws.Range("N14").Value = 3.25
$ws.Range("O14").Value = 1.39
$ws.Range("P14").Value = 1.74
$ws.Range("Q14").Value = 2.18
$ws.Range("X14").Value = 15
$ws.Range("AB14").Value = 9
$ws.Range("AD14").Value = 23
$ws.Range("AH14").Value = 27
$ws.Range("AL14").Value = 46

